# Holidays add to calendar - Done
# Duplicates the existing 10-row holiday block (rows 121-130) two more times,
# appending rows 131-140 and 141-150 to the worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The repeating 10-row block of data (columns A-F), taken from rows 121-130.
$blockData = @(
    @("Tr002", "2016-05-16", 8,   5, 8,   0),
    @("Tr003", "2016-05-17", 8.3, 5, 7.3, 0),
    @("Tr004", "2016-05-18", 8.4, 5, 7.2, 0),
    @("Tr005", "2016-05-19", 8.1, 5, 7.5, 0),
    @("Tr006", "2016-05-20", 8,   5, 8,   0),
    @("Tr007", "2016-05-21", 8.2, 5, 7.4, 0),
    @("Tr010", "2016-05-18", 8.4, 5, 7.2, 0),
    @("Tr011", "2016-05-19", 8.1, 5, 7.5, 0),
    @("Tr012", "2016-05-20", 8,   5, 8,   0),
    @("Tr013", "2016-05-21", 8.2, 5, 7.4, 0)
)

# Repeat the block twice (rows 131-140, then 141-150).
$startRow = 131
for ($rep = 0; $rep -lt 2; $rep++) {
    foreach ($item in $blockData) {
        $r = $startRow

        # Column A: transaction code (text)
        $ws.Cells.Item($r, 1).Value = $item[0]

        # Column B: date-looking text value - force text format first so it
        # isn't auto-converted into a date serial number by the Value setter,
        # then clear the explicit format so no style index is left on the cell.
        $ws.Cells.Item($r, 2).NumberFormat = "@"
        $ws.Cells.Item($r, 2).Value = $item[1]
        $ws.Cells.Item($r, 2).ClearFormats()

        # Columns C-F: numeric values
        $ws.Cells.Item($r, 3).Value = $item[2]
        $ws.Cells.Item($r, 4).Value = $item[3]
        $ws.Cells.Item($r, 5).Value = $item[4]
        $ws.Cells.Item($r, 6).Value = $item[5]

        $startRow = $startRow + 1
    }
}
